# "Fruta / hortaliza, semanal"
# Insert a new weekly record at the top of the Espinaca price history
# (row 247), pushing the existing rows 247-301 down to 248-302.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 247; this shifts every
# row from 247..301 down by one (to 248..302) and extends the sheet's
# used range to row 302.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new weekly data point.
$ws.Cells.Item(247, 1).Value = 8
$ws.Cells.Item(247, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(247, 3).Value = "Coquimbo"
$ws.Cells.Item(247, 4).Value = 44785
$ws.Cells.Item(247, 5).Value = 4
$ws.Cells.Item(247, 6).Value = 100112012
$ws.Cells.Item(247, 7).Value = "Espinaca"
$ws.Cells.Item(247, 8).Value = "Sin especificar"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 2000
$ws.Cells.Item(247, 11).Value = 500
$ws.Cells.Item(247, 12).Value = 550
$ws.Cells.Item(247, 13).Value = 525
$ws.Cells.Item(247, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(247, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(247, 16).Value = 1050
$ws.Cells.Item(247, 17).Value = 0.5
$ws.Cells.Item(247, 18).Value = "Hortaliza"

# Match the date-number format used by the rest of column D.
$ws.Cells.Item(247, 4).NumberFormat = $ws.Cells.Item(248, 4).NumberFormat
